$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.17"
$ws.Range("E2").Value = "'0.48%"
$ws.Range("D3").Value = "'36.28"
$ws.Range("E3").Value = "'0.88%"
$ws.Range("D4").Value = "'5.046"
$ws.Range("D5").Value = "'0.08150"
$ws.Range("E5").Value = "'0.68%"
$ws.Range("D6").Value = "'2.011"
$ws.Range("E6").Value = "'4.73%"
$ws.Range("D7").Value = "'4.163"
$ws.Range("E7").Value = "'0.37%"
$ws.Range("D8").Value = "'7.876"
$ws.Range("E8").Value = "'0.15%"
$ws.Range("D9").Value = "'0.9284"
$ws.Range("E9").Value = "'-0.21%"
$ws.Range("D10").Value = "'0.1480"
$ws.Range("E10").Value = "'18.17%"
$ws.Range("D11").Value = "'0.1932"
$ws.Range("E11").Value = "'1.12%"
$ws.Range("D12").Value = "'0.09109"
$ws.Range("E12").Value = "'-1.58%"
$ws.Range("D13").Value = "'0.03448"
$ws.Range("E13").Value = "'-1.60%"
$ws.Range("D14").Value = "'0.09885"
$ws.Range("E14").Value = "'-0.37%"
$ws.Range("D15").Value = "'0.001406"
$ws.Range("E15").Value = "'-1.27%"
$ws.Range("D16").Value = "'0.006467"
$ws.Range("E16").Value = "'1.55%"
$ws.Range("D17").Value = "'3.842"
$ws.Range("E17").Value = "'6.47%"
$ws.Range("D18").Value = "'3.401"
$ws.Range("E18").Value = "'5.65%"
$ws.Range("D19").Value = "'0.3465"
$ws.Range("E19").Value = "'0.77%"
$ws.Range("D20").Value = "'0.1320"
$ws.Range("E20").Value = "'2.02%"
$ws.Range("D21").Value = "'4.810"
$ws.Range("E21").Value = "'-7.28%"
$ws.Range("D22").Value = "'0.2338"
$ws.Range("E22").Value = "'-7.62%"
$ws.Range("D23").Value = "'0.04382"
$ws.Range("E23").Value = "'-0.66%"
$ws.Range("D24").Value = "'0.001232"
$ws.Range("E24").Value = "'-0.10%"
$ws.Range("E25").Value = "'-11.27%"
$ws.Range("D27").Value = "'0.0001298"
$ws.Range("E27").Value = "'0.00%"
$ws.Range("D39").Value = "'0.02053"
$ws.Range("E39").Value = "'4.62%"
$ws.Range("D40").Value = "'0.05172"
$ws.Range("E40").Value = "'-1.39%"
$ws.Range("D41").Value = "'0.007466"
$ws.Range("E41").Value = "'-1.22%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'-0.68%"
$ws.Range("E43").Value = "'0.30%"
$ws.Range("D44").Value = "'0.002127"
$ws.Range("E44").Value = "'1.42%"
$ws.Range("D45").Value = "'0.009676"
$ws.Range("E45").Value = "'-9.40%"
$ws.Range("D46").Value = "'0.00006293"
$ws.Range("E46").Value = "'-1.09%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'63.82"
$ws.Range("E48").Value = "'0.39%"
$ws.Range("D49").Value = "'0.001598"
$ws.Range("E49").Value = "'-3.62%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.08%"
